$wb = $excel.ActiveWorkbook

# --- Overview sheet: Latest HO Xliff Generate Date (column G) for rows 4-7 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-08-31 14:39:25"

# --- zh-cn sheet: Priority (column E) low -> ht, Latest Handoff Datetime (column H) for rows 4-7 ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4:E7").Value = "ht"
$wsZh.Range("H4:H7").Value = "2016-08-31 14:39:20"

# --- de-de sheet: Priority (column E) low -> ht for rows 4-7 ---
# Also Latest Handoff Datetime (column H) shares the same underlying string
# as the Overview sheet's G column for these rows, so keep them in sync.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4:E7").Value = "ht"
$wsDe.Range("H4:H7").Value = "2016-08-31 14:39:25"
